$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 35.00391733333333
$ws.Range("H2").Value = 105.011752
$ws.Range("I2").Value = 0.9591895364534718
$ws.Range("J2").Value = 0.9591895364534718
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.07605
$ws.Range("N2").Value = 18.22815
$ws.Range("O2").Value = 0.0302610603580868
$ws.Range("P2").Value = 0.0302610603580868
$ws.Range("Q2").Value = 212.6855519132
$ws.Range("R2").Value = 1914.1699672188
$ws.Range("S2").Value = 0.02902609245746381
$ws.Range("T2").Value = 0.02902609245746381

$ws.Range("G3").Value = 35.00391733333333
$ws.Range("H3").Value = 105.011752
$ws.Range("I3").Value = 0.9591895364534718
$ws.Range("J3").Value = 0.9591895364534718
$ws.Range("O3").Value = 0.2994824511432495
$ws.Range("P3").Value = 0.2994824511432494
$ws.Range("Q3").Value = 2104.869745342492
$ws.Range("R3").Value = 18943.82770808243
$ws.Range("S3").Value = 0.287260433488043
$ws.Range("T3").Value = 0.2872604334880429

$ws.Range("G4").Value = 35.00391733333333
$ws.Range("H4").Value = 105.011752
$ws.Range("I4").Value = 0.9591895364534718
$ws.Range("J4").Value = 0.9591895364534718
$ws.Range("M4").Value = 134.5792873333333
$ws.Range("N4").Value = 403.737862
$ws.Range("O4").Value = 0.6702564884986638
$ws.Range("P4").Value = 0.6702564884986637
$ws.Range("Q4").Value = 4710.802248594913
$ws.Range("R4").Value = 42397.22023735422
$ws.Range("S4").Value = 0.6429030105079652
$ws.Range("T4").Value = 0.6429030105079651

$ws.Range("I5").Value = 0.0008369499257158872
$ws.Range("J5").Value = 0.0008369499257158872
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.07605
$ws.Range("N5").Value = 18.22815
$ws.Range("O5").Value = 0.0302610603580868
$ws.Range("P5").Value = 0.0302610603580868
$ws.Range("Q5").Value = 0.18558079515
$ws.Range("R5").Value = 1.67022715635
$ws.Range("S5").Value = [double]"2.532699221878473E-05"
$ws.Range("T5").Value = [double]"2.532699221878472E-05"

$ws.Range("I6").Value = 0.0008369499257158872
$ws.Range("J6").Value = 0.0008369499257158872
$ws.Range("O6").Value = 0.2994824511432495
$ws.Range("P6").Value = 0.2994824511432494
$ws.Range("S6").Value = 0.0002506518152375544
$ws.Range("T6").Value = 0.0002506518152375544

$ws.Range("I7").Value = 0.0008369499257158872
$ws.Range("J7").Value = 0.0008369499257158872
$ws.Range("M7").Value = 134.5792873333333
$ws.Range("N7").Value = 403.737862
$ws.Range("O7").Value = 0.6702564884986638
$ws.Range("P7").Value = 0.6702564884986637
$ws.Range("Q7").Value = 4.110455173022
$ws.Range("R7").Value = 36.994096557198
$ws.Range("S7").Value = 0.000560971118259548
$ws.Range("T7").Value = 0.000560971118259548

$ws.Range("G8").Value = 1.458762333333333
$ws.Range("H8").Value = 4.376287
$ws.Range("I8").Value = 0.03997351362081222
$ws.Range("J8").Value = 0.03997351362081222
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.07605
$ws.Range("N8").Value = 18.22815
$ws.Range("O8").Value = 0.0302610603580868
$ws.Range("P8").Value = 0.0302610603580868
$ws.Range("Q8").Value = 8.863512875449999
$ws.Range("R8").Value = 79.77161587904999
$ws.Range("S8").Value = 0.001209640908404203
$ws.Range("T8").Value = 0.001209640908404203

$ws.Range("G9").Value = 1.458762333333333
$ws.Range("H9").Value = 4.376287
$ws.Range("I9").Value = 0.03997351362081222
$ws.Range("J9").Value = 0.03997351362081222
$ws.Range("O9").Value = 0.2994824511432495
$ws.Range("P9").Value = 0.2994824511432494
$ws.Range("Q9").Value = 87.71888791299911
$ws.Range("R9").Value = 789.469991216992
$ws.Range("S9").Value = 0.01197136583996891
$ws.Range("T9").Value = 0.01197136583996891

$ws.Range("G10").Value = 1.458762333333333
$ws.Range("H10").Value = 4.376287
$ws.Range("I10").Value = 0.03997351362081222
$ws.Range("J10").Value = 0.03997351362081222
$ws.Range("M10").Value = 134.5792873333333
$ws.Range("N10").Value = 403.737862
$ws.Range("O10").Value = 0.6702564884986638
$ws.Range("P10").Value = 0.6702564884986637
$ws.Range("Q10").Value = 196.3191952087104
$ws.Range("R10").Value = 1766.872756878394
$ws.Range("S10").Value = 0.0267925068724391
$ws.Range("T10").Value = 0.0267925068724391
